# Update countries & provincias Spain
# Applies the 3-April-2020 04:50 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 04:50"

# --- Plain numeric refresh (no country reordering) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 245080
$ws.Range("C4").Value = 203
$ws.Range("D4").Value = 10403
$ws.Range("E4").Value = 228602
$ws.Range("F4").Value = 5421
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 6075

# Row 8: China
$ws.Range("B8").Value = 81620
$ws.Range("C8").Value = 31
$ws.Range("D8").Value = 76571
$ws.Range("E8").Value = 1727
$ws.Range("F8").Value = 379
$ws.Range("G8").Value = 4
$ws.Range("H8").Value = 3322

# Row 20: Brasil
$ws.Range("B20").Value = 8066
$ws.Range("C20").Value = 22
$ws.Range("D20").Value = 127
$ws.Range("E20").Value = 7612
$ws.Range("F20").Value = 296
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 327

# Row 44: Mexico
$ws.Range("B44").Value = 1510
$ws.Range("C44").Value = 132
$ws.Range("D44").Value = 633
$ws.Range("E44").Value = 827
$ws.Range("F44").Value = 1
$ws.Range("G44").Value = 13
$ws.Range("H44").Value = 50

# --- Rows whose ranking swapped with their neighbour: update country name + data ---

# Row 35/36: Japon overtakes India
$ws.Range("A35").Value = "Japon"
$ws.Range("B35").Value = 2617
$ws.Range("C35").Value = 122
$ws.Range("D35").Value = 472
$ws.Range("E35").Value = 2082
$ws.Range("F35").Value = 60
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 63

$ws.Range("A36").Value = "India"
$ws.Range("B36").Value = 2543
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 191
$ws.Range("E36").Value = 2280
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 72

# Row 119/120: Camboya overtakes Kenia
$ws.Range("A119").Value = "Camboya"
$ws.Range("B119").Value = 114
$ws.Range("C119").Value = 4
$ws.Range("D119").Value = 35
$ws.Range("E119").Value = 79
$ws.Range("F119").Value = 1
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 0

$ws.Range("A120").Value = "Kenia"
$ws.Range("B120").Value = 110
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 4
$ws.Range("E120").Value = 103
$ws.Range("F120").Value = 2
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 3

# Row 122/123: Trinidad yTobago overtakes Isla de Man
$ws.Range("A122").Value = "Trinidad yTobago"
$ws.Range("B122").Value = 97
$ws.Range("C122").Value = 3
$ws.Range("D122").Value = 1
$ws.Range("E122").Value = 90
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 6

$ws.Range("A123").Value = "Isla de Man"
$ws.Range("B123").Value = 95
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 0
$ws.Range("E123").Value = 94
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 1

# Row 137/138: El Salvador overtakes Uganda
$ws.Range("A137").Value = "El Salvador"
$ws.Range("B137").Value = 46
$ws.Range("C137").Value = 5
$ws.Range("D137").Value = 0
$ws.Range("E137").Value = 44
$ws.Range("F137").Value = 4
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 2

$ws.Range("A138").Value = "Uganda"
$ws.Range("B138").Value = 45
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 0
$ws.Range("E138").Value = 45
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 0
